$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "عدد الحصص"

for ($r = 2; $r -le 35; $r++) {
    $ws.Cells.Item($r, 5).Value = 2
}

$ws.Range("E2").Select()
